$wb = $excel.ActiveWorkbook

# --- Hoja1!A1: update the conversion-of-the-day message with new rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 15.15 = 64001.52 pesos`n✅ 64001.52 pesos = 15.0 = 959.73 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas sheet: updated rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 66
$wsTasas.Range("O10").Value = 4224.1
$wsTasas.Range("N12").Value = 4267.98
$wsTasas.Range("O12").Value = 64
